# RRPONTSYD.xlsx data refresh (FRED series RRPONTSYD) + SeriesInfo metadata bump.
#
# "Data" sheet: append the new daily observations (2023-10-30 .. 2023-11-15,
# business days only) that FRED published after the prior snapshot, continuing
# straight on from the existing last row (518).
#
# "SeriesInfo" sheet: refresh the realtime_start/realtime_end/observation_end/
# last_updated/popularity metadata to match the new pull.

$wb = $excel.ActiveWorkbook
$wsData = $wb.Worksheets.Item("Data")

# Seed rows 519:531 with the same style the rest of column A's date cells use
# (numeric date format + bold/centered header-ish look carried down the
# column) by copying the last existing row's formatting onto the new block
# before any values are written into it.
$wsData.Range("A518").Copy($wsData.Range("A519:A531"))

$wsData.Range("A519").Value = 45229
$wsData.Range("B519").Value = 1138.035
$wsData.Range("A520").Value = 45230
$wsData.Range("B520").Value = 1137.697
$wsData.Range("A521").Value = 45231
$wsData.Range("B521").Value = 1079.462
$wsData.Range("A522").Value = 45232
$wsData.Range("B522").Value = 1054.986
$wsData.Range("A523").Value = 45233
$wsData.Range("B523").Value = 1071.139
$wsData.Range("A524").Value = 45236
$wsData.Range("B524").Value = 1062.878
$wsData.Range("A525").Value = 45237
$wsData.Range("B525").Value = 1008.685
$wsData.Range("A526").Value = 45238
$wsData.Range("B526").Value = 1024.451
$wsData.Range("A527").Value = 45239
$wsData.Range("B527").Value = 993.314
$wsData.Range("A528").Value = 45240
$wsData.Range("B528").Value = 1032.72
$wsData.Range("A529").Value = 45243
$wsData.Range("B529").Value = 1020.272
$wsData.Range("A530").Value = 45244
$wsData.Range("B530").Value = 988.298
$wsData.Range("A531").Value = 45245
$wsData.Range("B531").Value = 944.241

$wsInfo = $wb.Worksheets.Item("SeriesInfo")

# B3 / B4 / B7 hold plain "YYYY-MM-DD" strings. Assigning that text straight
# to .Value would let Excel's own auto-detection turn it into a date serial
# (and swap in a date number format), so format the cell as Text first to
# keep it a literal string, then clear the format back off again so the
# cell's style is left exactly as it was before (unstyled), matching the
# source edit which only changes the displayed text.
$wsInfo.Range("B3").NumberFormat = "@"
$wsInfo.Range("B3").Value = "2023-11-15"
$wsInfo.Range("B3").ClearFormats()

$wsInfo.Range("B4").NumberFormat = "@"
$wsInfo.Range("B4").Value = "2023-11-15"
$wsInfo.Range("B4").ClearFormats()

$wsInfo.Range("B7").NumberFormat = "@"
$wsInfo.Range("B7").Value = "2023-11-15"
$wsInfo.Range("B7").ClearFormats()

# B14's value carries a UTC-offset suffix, which Excel does not parse as a
# recognised date/time literal, so it is safe to assign directly.
$wsInfo.Range("B14").Value = "2023-11-15 13:01:02-06"

# B15 (popularity) is a plain number.
$wsInfo.Range("B15").Value = 93
